$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A28").Value = "debut_couteaux (mm)"
$ws.Range("B28").Value = 10

$ws.Range("A29").Value = "r_int_BAR (mm)"
$ws.Range("B29").Value = 3

$ws.Range("A30").Value = "r_ext_BAR (mm)"
$ws.Range("B30").Value = 4

$ws.Range("A31").Value = "r_ini_couteaux (mm)"
$ws.Range("B31").Value = 3

$ws.Range("A28:A31").NumberFormat = $ws.Range("A27").NumberFormat
$ws.Range("B28:B31").NumberFormat = $ws.Range("B27").NumberFormat

$ws.Range("C32").Select()
